$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "hoorceollege" (oefeningen tijd) hours for week 2 (D) and week 3 (E)
$ws.Range("D12").Value = 0.10416666666666667
$ws.Range("D13").Value = 0.10416666666666667
$ws.Range("D14").Value = 0.10416666666666667
$ws.Range("D15").Value = 0.10416666666666667

$ws.Range("E12").Value = 0.16666666666666666
$ws.Range("E13").Value = 0.16666666666666666
$ws.Range("E14").Value = 0.16666666666666666
$ws.Range("E15").Value = 0.10416666666666667

# Update the selected cell to reflect the new active cell in the saved view
$ws.Range("F15").Select()
